$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1961
$ws.Range("I15").Value = 1961
$ws.Range("K15").Value = 5883
$ws.Range("M15").Value = -5714
$ws.Range("H40").Value = 8071.6665
$ws.Range("I40").Value = 6375
$ws.Range("J40").Value = 9429
$ws.Range("K40").Value = 6375
$ws.Range("L40").Value = 9429
$ws.Range("M40").Value = -6200
$ws.Range("N40").Value = -9779
$ws.Range("H93").Value = 25000
$ws.Range("J93").Value = 25000
$ws.Range("L93").Value = 25000
$ws.Range("N93").Value = -29992
$ws.Range("H96").Value = 888.2857
$ws.Range("I96").Value = 773.8
$ws.Range("K96").Value = 2321.4
$ws.Range("M96").Value = -948.3999999999996
$ws.Range("H111").Value = 1838.6945
$ws.Range("I111").Value = 6049.6665
$ws.Range("J111").Value = 996.5
$ws.Range("K111").Value = 18148.9995
$ws.Range("L111").Value = 2989.5
$ws.Range("M111").Value = -15081.9995
$ws.Range("N111").Value = -9123.5
$ws.Range("H116").Value = 3723.5
$ws.Range("I116").Value = 3581.6667
$ws.Range("K116").Value = 3581.6667
$ws.Range("M116").Value = -139.6667000000002
$ws.Range("H136").Value = 94354.5
$ws.Range("J136").Value = 103000
$ws.Range("L136").Value = 103000
$ws.Range("N136").Value = -113200
$ws.Range("H138").Value = 2520.3076
$ws.Range("I138").Value = 2193
$ws.Range("K138").Value = 6579
$ws.Range("M138").Value = -1439
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5776.778
$ws.Range("I61").Value = 2997.75
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 2997.75
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -2785.75
$ws.Range("N61").Value = -8424
$ws.Range("H110").Value = 2870.4167
$ws.Range("I110").Value = 2044.4
$ws.Range("J110").Value = 7000.5
$ws.Range("K110").Value = 2044.4
$ws.Range("L110").Value = 7000.5
$ws.Range("M110").Value = 0.5999999999999091
$ws.Range("N110").Value = -11090.5
$ws.Range("H122").Value = 3300
$ws.Range("I122").Value = 2950
$ws.Range("K122").Value = 8850
$ws.Range("M122").Value = -6400
$ws.Range("H132").Value = 1999.5
$ws.Range("I132").Value = 999
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 2997
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -467
$ws.Range("N132").Value = -14060
$ws.Range("H135").Value = 85942
$ws.Range("J135").Value = 91997
$ws.Range("L135").Value = 91997
$ws.Range("N135").Value = -102137
$ws.Range("H136").Value = 5776.778
$ws.Range("I136").Value = 2997.75
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 8993.25
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -6443.25
$ws.Range("N136").Value = -29100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2988.3157
$ws.Range("I86").Value = 2805.2307
$ws.Range("J86").Value = 3385
$ws.Range("K86").Value = 2805.2307
$ws.Range("L86").Value = 3385
$ws.Range("M86").Value = -1682.2307
$ws.Range("N86").Value = -5631
$ws.Range("H89").Value = 2988.3157
$ws.Range("I89").Value = 2805.2307
$ws.Range("J89").Value = 3385
$ws.Range("K89").Value = 14026.1535
$ws.Range("L89").Value = 16925
$ws.Range("M89").Value = -8410.1535
$ws.Range("N89").Value = -28157
$ws.Range("H134").Value = 4671
$ws.Range("I134").Value = 1999.5
$ws.Range("J134").Value = 10014
$ws.Range("K134").Value = 5998.5
$ws.Range("L134").Value = 30042
$ws.Range("M134").Value = -3463.5
$ws.Range("N134").Value = -35112
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1026
$ws.Range("I5").Value = 1026
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1026
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -914
$ws.Range("N5").ClearContents()
$ws.Range("H31").Value = 4258.7354
$ws.Range("I31").Value = 2626.8096
$ws.Range("K31").Value = 2626.8096
$ws.Range("M31").Value = -2331.8096
$ws.Range("H34").Value = 4258.7354
$ws.Range("I34").Value = 2626.8096
$ws.Range("K34").Value = 2626.8096
$ws.Range("M34").Value = -2424.8096
$ws.Range("H99").Value = 4055.2727
$ws.Range("I99").Value = 3822.6
$ws.Range("K99").Value = 3822.6
$ws.Range("M99").Value = -2324.6
$ws.Range("H126").Value = 4055.2727
$ws.Range("I126").Value = 3822.6
$ws.Range("K126").Value = 11467.8
$ws.Range("M126").Value = -8997.799999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 37129330
$ws.Range("J9").Value = 1550
$ws.Range("L9").Value = 4650
$ws.Range("N9").Value = -5098
$ws.Range("H23").Value = 6896.5
$ws.Range("I23").Value = 436.66666
$ws.Range("J23").Value = 9049.777
$ws.Range("K23").Value = 1309.99998
$ws.Range("L23").Value = 27149.331
$ws.Range("M23").Value = -1074.99998
$ws.Range("N23").Value = -27619.331
$ws.Range("H32").Value = 97225820
$ws.Range("J32").Value = 15155447
$ws.Range("L32").Value = 45466341
$ws.Range("N32").Value = -45466907
$ws.Range("H140").Value = 12227.25
$ws.Range("I140").Value = 5563.6
$ws.Range("K140").Value = 16690.8
$ws.Range("M140").Value = -11510.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 45456984
$ws.Range("I80").Value = 100002330
$ws.Range("J80").Value = 2526.75
$ws.Range("K80").Value = 100002330
$ws.Range("L80").Value = 2526.75
$ws.Range("M80").Value = -100001332
$ws.Range("N80").Value = -4522.75
$ws.Range("H83").Value = 45456984
$ws.Range("I83").Value = 100002330
$ws.Range("J83").Value = 2526.75
$ws.Range("K83").Value = 500011650
$ws.Range("L83").Value = 12633.75
$ws.Range("M83").Value = -500006658
$ws.Range("N83").Value = -22617.75
$ws.Range("H96").Value = 39387
$ws.Range("J96").Value = 39387
$ws.Range("L96").Value = 39387
$ws.Range("N96").Value = -44879
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2199.6667
$ws.Range("I16").Value = 443
$ws.Range("J16").Value = 3078
$ws.Range("K16").Value = 443
$ws.Range("L16").Value = 3078
$ws.Range("M16").Value = -273
$ws.Range("N16").Value = -3418
$ws.Range("H40").Value = 70001
$ws.Range("I40").Value = 70001
$ws.Range("K40").Value = 70001
$ws.Range("M40").Value = -69865
$ws.Range("H61").Value = 3145.3215
$ws.Range("I61").Value = 3015.8696
$ws.Range("J61").Value = 3740.8
$ws.Range("K61").Value = 3015.8696
$ws.Range("L61").Value = 3740.8
$ws.Range("M61").Value = -2813.8696
$ws.Range("N61").Value = -4144.8
$ws.Range("H113").Value = 3145.3215
$ws.Range("I113").Value = 3015.8696
$ws.Range("J113").Value = 3740.8
$ws.Range("K113").Value = 3015.8696
$ws.Range("L113").Value = 3740.8
$ws.Range("M113").Value = -845.8696
$ws.Range("N113").Value = -8080.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3690.6316
$ws.Range("I81").Value = 2718.5833
$ws.Range("J81").Value = 5357
$ws.Range("K81").Value = 5437.1666
$ws.Range("L81").Value = 10714
$ws.Range("M81").Value = -4376.1666
$ws.Range("N81").Value = -12836
$ws.Range("H84").Value = 3690.6316
$ws.Range("I84").Value = 2718.5833
$ws.Range("J84").Value = 5357
$ws.Range("K84").Value = 27185.833
$ws.Range("L84").Value = 53570
$ws.Range("M84").Value = -21881.833
$ws.Range("N84").Value = -64178
$ws.Range("H132").Value = 2477.5264
$ws.Range("I132").Value = 2135.8125
$ws.Range("J132").Value = 4300
$ws.Range("K132").Value = 6407.4375
$ws.Range("L132").Value = 12900
$ws.Range("M132").Value = -3877.4375
$ws.Range("N132").Value = -17960

Write-Output "Applied 201 cell edits"
